# Fix bug in fedrollover: append the missing 2009 rows (Jan-Jul) that were
# dropped off the end of the SOMA series on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 380
$startDate = 20090100
$step = 100
$count = 7

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $dateVal = $startDate + ($i * $step)
    $ws.Cells.Item($row, 1).Value = $dateVal
    $ws.Cells.Item($row, 2).Value = 0
}
